$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.605.20"
$ws.Range("E2").Value = "  -1.55%  "
$ws.Range("D3").Value = "3.025.45"
$ws.Range("E3").Value = "  -0.57%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.998"
$ws.Range("E4").Value = "  -0.20%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "533.35"
$ws.Range("E5").Value = "  +0.05%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "134.12"
$ws.Range("E6").Value = "  +1.67%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("D8").Value = "3.027.20"
$ws.Range("E8").Value = "  -0.20%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.495"
$ws.Range("E9").Value = "  +0.13%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.150"
$ws.Range("E10").Value = "  -1.78%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.19"
$ws.Range("E11").Value = "  +1.17%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.447"
$ws.Range("E12").Value = "  -0.54%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000221"
$ws.Range("E13").Value = "  -0.09%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.27"
$ws.Range("E14").Value = "  +1.22%  "
$ws.Range("D15").Value = "3.499.72"
$ws.Range("E15").Value = "  -0.92%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.111"
$ws.Range("E16").Value = "  +0.23%  "
$ws.Range("D17").Value = "61.498.42"
$ws.Range("E17").Value = "  -1.62%  "
$ws.Range("D18").Value = "3.002.97"
$ws.Range("E18").Value = "  -1.24%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.62"
$ws.Range("E19").Value = "  +0.99%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "467.05"
$ws.Range("E20").Value = "  -2.90%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.30"
$ws.Range("E21").Value = "  +1.23%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.679"
$ws.Range("E22").Value = "  -0.98%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.96"
$ws.Range("E23").Value = "  -1.21%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "79.35"
$ws.Range("E24").Value = "  +0.79%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.04"
$ws.Range("E25").Value = "  +0.33%  "
$ws.Range("E26").Value = "  +0.44%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.68"
$ws.Range("E27").Value = "  +0.11%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.87"
$ws.Range("E28").Value = "  -1.55%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.995"
$ws.Range("E29").Value = "  -0.36%  "
$ws.Range("B30").Value = "ImmutableX"
$ws.Range("C30").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.89"
$ws.Range("E30").Value = "  +2.73%  "
$ws.Range("B31").Value = "EthereumClassic"
$ws.Range("C31").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "25.65"
$ws.Range("E31").Value = "  -0.39%  "
$ws.Range("B32").Value = "Mantle"
$ws.Range("C32").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.15"
$ws.Range("E32").Value = "  +4.05%  "
$ws.Range("B33").Value = "NEARProtocol"
$ws.Range("C33").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.51"
$ws.Range("E33").Value = "  +4.50%  "
$ws.Range("B34").Value = "OKB"
$ws.Range("C34").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "55.62"
$ws.Range("E34").Value = "  -1.53%  "
$ws.Range("B35").Value = "Stacks"
$ws.Range("C35").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.29"
$ws.Range("E35").Value = "  -1.90%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.90"
$ws.Range("E36").Value = "  +0.17%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "464.83"
$ws.Range("E37").Value = "  -1.49%  "
$ws.Range("D38").Value = "3.207.03"
$ws.Range("E38").Value = "  +4.71%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0790"
$ws.Range("E39").Value = "  +0.13%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0386"
$ws.Range("E40").Value = "  -1.31%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.117"
$ws.Range("E41").Value = "  +2.97%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.15"
$ws.Range("E42").Value = "  +1.78%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.49"
$ws.Range("E43").Value = "  -4.63%  "
$ws.Range("B44").Value = "InjectiveProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "26.78"
$ws.Range("E44").Value = "  +10.67%  "
$ws.Range("B45").Value = "USDe"
$ws.Range("C45").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.00"
$ws.Range("E45").Value = "  +0.05%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.248"
$ws.Range("E46").Value = "  -0.30%  "
$ws.Range("B47").Value = "Fetch.AI"
$ws.Range("C47").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.01"
$ws.Range("E47").Value = "  +0.50%  "
$ws.Range("B48").Value = "Monero"
$ws.Range("C48").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "119.76"
$ws.Range("E48").Value = "  -0.49%  "
$ws.Range("B49").Value = "Stellar"
$ws.Range("C49").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.109"
$ws.Range("E49").Value = "  +1.29%  "
$ws.Range("D50").Value = "0.0₃0503"
$ws.Range("E50").Value = "  -5.81%  "
$ws.Range("E51").Value = "  +8.50%  "
